# Update Excel files after daily scrape - 2025-09-16 03:04:48 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes (C, D, F, G, H) -------------------------------
# Stored OOXML "width" = ColumnWidth + 0.8333333333333334, so back the
# offset out to land on the exact target stored widths.
$colWidthOffset = 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 90 - $colWidthOffset
$ws.Columns.Item(4).ColumnWidth = 70 - $colWidthOffset
$ws.Columns.Item(6).ColumnWidth = 17 - $colWidthOffset
$ws.Columns.Item(7).ColumnWidth = 16 - $colWidthOffset
$ws.Columns.Item(8).ColumnWidth = 60 - $colWidthOffset

# --- Data rows -----------------------------------------------------------
# Column A holds numeric-looking opportunity IDs that must stay text, so
# force a text number-format before assigning, then reset the style back
# to "Normal" so no stray cell style/format is left behind.
$idRange = $ws.Range("A2:A11")
$idRange.NumberFormat = "@"

$data = @(
    @("1327205", "https://aiesec.org/opportunity/global-talent/1327205", "Environmental Technician", "Calgary, AB, Canada", "No", "10 applicants", "6 - 18 Months", "Oak Environmental Inc."),
    @("1326923", "https://aiesec.org/opportunity/global-talent/1326923", "AI Tech Developer", "Panamá, Provincia de Panamá, Panamá", "No", "62 applicants", "6 - 18 Months", "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"),
    @("1326590", "https://aiesec.org/opportunity/global-talent/1326590", "BUSINESS DEVELOPMENT", "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt", "No", "24 applicants", "9 - 12 Weeks", "Egypt holiday travel"),
    @("1326575", "https://aiesec.org/opportunity/global-talent/1326575", "Graphic design", "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt", "No", "9 applicants", "9 - 12 Weeks", "Egypt holiday travel"),
    @("1326555", "https://aiesec.org/opportunity/global-talent/1326555", "Graphic desgin", "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt", "No", "1 applicant", "9 - 12 Weeks", "Egypt holiday travel"),
    @("1326536", "https://aiesec.org/opportunity/global-talent/1326536", "Marketing", "Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt", "No", "6 applicants", "9 - 12 Weeks", "Egypt holiday travel"),
    @("1326357", "https://aiesec.org/opportunity/global-talent/1326357", "Human resource strategies to attract, engage, and inspire through the lens of marketing", "日本、京都府京都市", "No", "18 applicants", "9 - 12 Weeks", "Shinko Automotive Co., Ltd."),
    @("1317404", "https://aiesec.org/opportunity/global-talent/1317404", "Demand Generation", "Mysuru, Karnataka, India", "No", "11 applicants", "6 - 18 Months", "DeUS Tech Services"),
    @("1316788", "https://aiesec.org/opportunity/global-talent/1316788", "Travel Coordinator", "Mexico City, CDMX, Mexico", "No", "103 applicants", "6 - 18 Months", "Ikan Experience"),
    @("1303327", "https://aiesec.org/opportunity/global-talent/1303327", "Administrative Assistant", "Konak, Türkiye", "No", "562 applicants", "9 - 12 Weeks", "Tekinalp Holding")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $rowIndex++
}

# Reset column A's style back to the default so no stray "Normal"-with-
# text-format style lingers on the written cells.
$idRange.Style = "Normal"
